$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat'
$ws.Range('G3').Value = 'Dr. Alshimaa Atef, Administrator, Dr. Manar Montaser, Dr. Gehan Adel'
$ws.Range('G4').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Hanan Ragab'
$ws.Range('G7').Value = 'Dr. Safa Hany, Dr. Amal Awwad'
$ws.Range('G10').Value = 'Dr. Basma Hamed, Dr. Amira Ibrahim'
$ws.Range('G12').Value = 'Dr. Dalia Tarek Elwan, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad'
$ws.Range('G17').Value = 'Dr. Enas Omran, Dr. Walaa Ghanima'
$ws.Range('G18').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Abdullah El-Agrody'
$ws.Range('G19').Value = 'Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Eman Samir Gabry'
$ws.Range('G20').Value = 'Dr. Nardine, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Monica, Dr. Marina Sorial, Dr. Remon, Dr. Marina Atef'
$ws.Range('G21').Value = 'Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat'
$ws.Range('G22').Value = 'Dr. Alshimaa Atef, Administrator, Dr. Manar Montaser, Dr. Gehan Adel'
$ws.Range('G23').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Hanan Ragab'
$ws.Range('G24').Value = 'Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Nada Mohammad'
$ws.Range('G25').Value = 'Dr. Yasmin Tarek, Dr. Nourhan Mohammad'
$ws.Range('G26').Value = 'Dr. Safa Hany, Dr. Amal Awwad'
$ws.Range('G28').Value = 'Dr. Eman M. Abo-Sakaya, Dr. Arwa Al-Sayed, Dr. Marwa Mustafa, Dr. Esraa Mostafa, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Dina Adel'
$ws.Range('G29').Value = 'Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Amira Ibrahim'
$ws.Range('G31').Value = 'Dr. Dalia Tarek Elwan, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad'
$ws.Range('G36').Value = 'Dr. Enas Omran, Dr. Walaa Ghanima'
$ws.Range('G37').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Abdullah El-Agrody'
$ws.Range('G38').Value = 'Dr. Nardine, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Monica, Dr. Marina Sorial, Dr. Remon, Dr. Marina Atef'
$ws.Range('G39').Value = 'Dr. Nardine, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Monica, Dr. Marina Sorial, Dr. Remon, Dr. Marina Atef'
$ws.Range('G40').Value = 'Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat'
$ws.Range('G41').Value = 'Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Hanan Ragab'
$ws.Range('G42').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi'
$ws.Range('G43').Value = 'Dr. Lamiaa Ossama, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Nada Mohammad'
$ws.Range('G45').Value = 'Dr. Safa Hany, Dr. Amal Awwad'
$ws.Range('G47').Value = 'Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Amira Ibrahim, Dr. Maryam Ahmad'
$ws.Range('G48').Value = 'Dr. Eman M. Abo-Sakaya, Dr. Merna Said, Dr. Amany Raafat, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Maryam Ahmad'
$ws.Range('G50').Value = 'Dr. Dalia Tarek Elwan, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad'
$ws.Range('G56').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Abdullah El-Agrody'
$ws.Range('G57').Value = 'Dr. Nardine, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Monica, Dr. Marina Sorial, Dr. Remon, Dr. Marina Atef'
$ws.Range('G58').Value = 'Dr. Nardine, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Monica, Dr. Marina Sorial, Dr. Remon, Dr. Marina Atef'
$ws.Range('G59').Value = 'Dr. Nourhan Mahmoud, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda'
$ws.Range('G60').Value = 'Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Hanan Ragab'
$ws.Range('G61').Value = 'Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Nahla Nagiub'
$ws.Range('G66').Value = 'Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Dina Adel'
$ws.Range('G67').Value = 'Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Amira Ibrahim'
$ws.Range('G71').Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range('G72').Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range('G75').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Abdullah El-Agrody'
$ws.Range('G76').Value = 'Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Eman Samir Gabry'
$ws.Range('G77').Value = 'Dr. Nardine, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Monica, Dr. Marina Sorial, Dr. Remon, Dr. Marina Atef'
$ws.Range('G78').Value = 'Dr. Nourhan Mahmoud, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda'
$ws.Range('G79').Value = 'Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Hanan Ragab'
$ws.Range('G80').Value = 'Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Nahla Nagiub'
$ws.Range('G81').Value = 'Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Nada Mohammad'
$ws.Range('G82').Value = 'Dr. Yasmin Tarek, Dr. Nourhan Mohammad'
$ws.Range('G83').Value = 'Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed'
$ws.Range('G85').Value = 'Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Dina Adel'
$ws.Range('G86').Value = 'Dr. Eman M. Abo-Sakaya, Dr. Merna Said, Dr. Amany Raafat, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Maryam Ahmad'
$ws.Range('G88').Value = 'Dr. Dalia Tarek Elwan, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad'
$ws.Range('G90').Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range('G91').Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range('G94').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Abdullah El-Agrody'
$ws.Range('G95').Value = 'Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Eman Samir Gabry'
$ws.Range('G96').Value = 'Dr. Nardine, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Monica, Dr. Marina Sorial, Dr. Remon, Dr. Marina Atef'
$ws.Range('G97').Value = 'Dr. Nourhan Mahmoud, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda'
$ws.Range('G98').Value = 'Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Hanan Ragab'
$ws.Range('G99').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi'
$ws.Range('G100').Value = 'Dr. Lamiaa Ossama, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Nada Mohammad'
$ws.Range('G102').Value = 'Dr. Safa Hany, Dr. Amal Awwad'
$ws.Range('G104').Value = 'Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Amira Ibrahim, Dr. Maryam Ahmad'
$ws.Range('G105').Value = 'Dr. Basma Hamed, Dr. Amira Ibrahim'
$ws.Range('G113').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Abdullah El-Agrody'
$ws.Range('G115').Value = 'Dr. Nardine, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Monica, Dr. Marina Sorial, Dr. Remon, Dr. Marina Atef'
